# Added new test case for Negative login

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet from "Login_positive" to "Login"
$ws.Name = "Login"

# Add the new negative-login test data row (row 3)
$ws.Cells.Item(3, 1).Value = "test@test.com"
$ws.Cells.Item(3, 2).Value = "test123"

# Turn the new username cell into a mailto hyperlink (adds the built-in
# "Hyperlink" cell style / font automatically)
$ws.Hyperlinks.Add($ws.Cells.Item(3, 1), "mailto:test@test.com") | Out-Null

# Match the final selection left behind in the sheet
$ws.Range("I13").Select() | Out-Null
